$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105 (shifts existing rows 105-206 down to 106-207)
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row 105 with the new price-report record
$ws.Range("A105").Value = 4
$ws.Range("B105").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C105").Value = "Los Lagos"
$ws.Range("D105").Value = 44789
$ws.Range("E105").Value = 10
$ws.Range("F105").Value = 100112009
$ws.Range("G105").Value = "Acelga"
$ws.Range("H105").Value = "Sin especificar"
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 220
$ws.Range("K105").Value = 1500
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 1500
$ws.Range("N105").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O105").Value = "Región de Los Lagos"
$ws.Range("P105").Value = 1000
$ws.Range("Q105").Value = 1.5
$ws.Range("R105").Value = "Hortaliza"
